$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 previously held the "useGravitySettling" T/F toggle; it was
# replaced with the new "gravity_species" input row (species-name :
# gravity coefficient, defaulted to SO4:gd(j)).
$ws.Range("A13").Value = "gravity_species"
$ws.Range("B13").Value = "species-name : gravity coeff (put # if no input)"
$ws.Range("C13").Value = "SO4:gd(j)"

# Reflect the author's final view state: zoomed out a bit and left the
# selection on D13.
$excel.ActiveWindow.Zoom = 145
$ws.Range("D13").Select() | Out-Null
